# Applies the changes described in the commit for A2_Structure_Lists.xlsx:
#  1. The "Timeslice" list value "All" (cell C2 on the "Lists" sheet) becomes "Some",
#     matching the xtra_scen Timeslice now used in MOMF_T1_A.yaml.
#  2. The "Emission" list (cells E2:E21 on the "Lists" sheet) is re-written with the
#     same 20 members but in a new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lists")

# 1) Timeslice list: "All" -> "Some"
$ws.Range("C2").Value = "Some"

# 2) Emission list: re-order the 20 existing entries (E2:E21)
$emissions = @(
    "CONTUR",
    "salud_residuos",
    "contam_agua",
    "RESHID",
    "CO2e_HFC",
    "CONHICK",
    "DAPANI",
    "FERT_ORG",
    "CO2e_PIUP",
    "CONVAR",
    "CO2e_sources",
    "CO2e_PP",
    "Health",
    "CO2e_DE",
    "turismo_residuos",
    "RM",
    "CO2e_TRN",
    "CONHAB",
    "CO2e_AFOLU",
    "CO2e_WASTE"
)

for ($i = 0; $i -lt $emissions.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $emissions[$i]
}
